$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.526080131530762
$ws.Range("B1").Value = 1.829545497894287
$ws.Range("C1").Value = 2.452124357223511
$ws.Range("D1").Value = 4.997652530670166
$ws.Range("E1").Value = 1.222767233848572
